$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nine")

$ws.Range("D3").Value = 10.69
$ws.Range("E3").Value = 10.19

$ws.Range("C4").Value = 9.31
$ws.Range("E4").Value = 9.94
$ws.Range("F4").Value = 10.09

$ws.Range("C5").Value = 9.81
$ws.Range("D5").Value = 10.06
$ws.Range("F5").Value = 10.22
$ws.Range("G5").Value = 9.24
$ws.Range("H5").Value = 7.88

$ws.Range("D6").Value = 9.91
$ws.Range("E6").Value = 9.779999999999999
$ws.Range("G6").Value = 10.15
$ws.Range("H6").Value = 11.96

$ws.Range("E7").Value = 10.76
$ws.Range("F7").Value = 9.85

$ws.Range("E8").Value = 12.12
$ws.Range("F8").Value = 8.039999999999999
